$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.354.72'
$ws.Range('E2').Value = '  -1.10%  '
$ws.Range('D3').Value = '1.802.84'
$ws.Range('E3').Value = '  -0.94%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '227.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('E6').Value = '  +4.12%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '36.37'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.40%  '
$ws.Range('E9').Value = '  -1.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0691'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.25%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '2.064.03'
$ws.Range('E12').Value = '  -0.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.31'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.23%  '
$ws.Range('D14').Value = '1.817.87'
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').Value = '34.367.64'
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.58'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.28%  '
$ws.Range('D20').Value = '0.0₃0786'
$ws.Range('E20').Value = '  -2.16%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.46'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('E23').Value = '  -0.57%  '
$ws.Range('E24').Value = '  +6.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '170.47'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '17.43'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.122'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.27%  '
$ws.Range('E29').Value = '  +0.23%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.98'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.24'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.45%  '
$ws.Range('E32').Value = '  -0.95%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0521'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.25%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.80'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.43%  '
$ws.Range('D35').Value = '1.379.01'
$ws.Range('E35').Value = '  -2.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.660'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.74%  '
$ws.Range('E37').Value = '  -1.44%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -11.44%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0187'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.55%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '82.19'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.62%  '
$ws.Range('B41').Value = 'ARBITRUM'
$ws.Range('C41').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.951'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.42'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.80'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.28%  '
$ws.Range('E44').Value = '  +6.67%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0498'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.90%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.94'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.76%  '
$ws.Range('D48').Value = '1.965.64'
$ws.Range('E48').Value = '  -0.88%  '
$ws.Range('E49').Value = '  +0.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.42%  '
$ws.Range('D51').Value = '0.0₆0124'
$ws.Range('E51').Value = '  -5.06%  '
